$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.719.13"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.581.78"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "602.11"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.75%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "137.54"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").Value = "3.581.76"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.59%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.18"
$c.ClearFormats()
$ws.Range("E11").Value = "  +5.20%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.392"
$c.ClearFormats()
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "4.190.24"
$ws.Range("E13").Value = "  +1.41%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.08"
$c.ClearFormats()
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "3.590.84"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "65.796.68"
$ws.Range("E18").Value = "  +0.81%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.01"
$c.ClearFormats()
$ws.Range("E19").Value = "  -2.88%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.61"
$c.ClearFormats()
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  -1.33%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "395.50"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("D24").Value = "3.726.84"
$ws.Range("E24").Value = "  +1.50%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "74.18"
$c.ClearFormats()
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  +2.84%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.12"
$c.ClearFormats()
$ws.Range("E28").Value = "  +5.21%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.64"
$c.ClearFormats()
$ws.Range("E29").Value = "  +30.47%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.ClearFormats()
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "3.585.44"
$ws.Range("E33").Value = "  +1.14%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "24.48"
$c.ClearFormats()
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  +2.02%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.38"
$c.ClearFormats()
$ws.Range("E37").Value = "  +7.85%  "
$ws.Range("E38").Value = "  +5.14%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "7.06"
$c.ClearFormats()
$ws.Range("E39").Value = "  +1.61%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "168.70"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +4.49%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.839"
$c.ClearFormats()
$ws.Range("E42").Value = "  +1.95%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "27.13"
$c.ClearFormats()
$ws.Range("E43").Value = "  +4.99%  "
$ws.Range("E44").Value = "  +8.62%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "43.11"
$c.ClearFormats()
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "2.451.54"
$ws.Range("E50").Value = "  +2.96%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.909"
$c.ClearFormats()
$ws.Range("E51").Value = "  +10.41%  "
